$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.939.23'
$ws.Range('E2').Value = '  +1.05%  '
$ws.Range('D3').Value = '2.358.07'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'0.692"
$ws.Range('E5').Value = '  +6.27%  '
$ws.Range('D6').Value = "'242.14"
$ws.Range('E6').Value = '  +3.20%  '
$ws.Range('D7').Value = "'77.02"
$ws.Range('E7').Value = '  +6.38%  '
$ws.Range('D9').Value = "'0.634"
$ws.Range('E9').Value = '  +26.79%  '
$ws.Range('E10').Value = '  +4.67%  '
$ws.Range('D11').Value = "'57.40"
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = "'33.94"
$ws.Range('E12').Value = '  +24.44%  '
$ws.Range('D13').Value = "'7.57"
$ws.Range('E13').Value = '  +19.79%  '
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').Value = '2.708.24'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('D16').Value = "'16.92"
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('D17').Value = "'0.925"
$ws.Range('E17').Value = '  +6.18%  '
$ws.Range('D18').Value = '2.354.86'
$ws.Range('E18').Value = '  -0.87%  '
$ws.Range('D19').Value = '43.881.96'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  +3.12%  '
$ws.Range('D21').Value = "'6.68"
$ws.Range('E21').Value = '  +4.68%  '
$ws.Range('D22').Value = "'77.72"
$ws.Range('E22').Value = '  +3.25%  '
$ws.Range('D23').Value = "'256.83"
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('D25').Value = "'2.55"
$ws.Range('E25').Value = '  +2.73%  '
$ws.Range('E26').Value = '  +10.21%  '
$ws.Range('D27').Value = "'3.62"
$ws.Range('E27').Value = '  -5.65%  '
$ws.Range('D28').Value = "'1.79"
$ws.Range('E28').Value = '  +15.83%  '
$ws.Range('E29').Value = '  +2.31%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('D31').Value = "'174.58"
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('E33').Value = '  +6.40%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'0.0758"
$ws.Range('E34').Value = '  +8.69%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = "'5.33"
$ws.Range('E35').Value = '  +5.90%  '
$ws.Range('D36').Value = "'5.41"
$ws.Range('E36').Value = '  +6.75%  '
$ws.Range('D37').Value = "'3.83"
$ws.Range('E37').Value = '  +2.14%  '
$ws.Range('D38').Value = "'2.44"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = "'6.49"
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  +8.25%  '
$ws.Range('D41').Value = "'19.35"
$ws.Range('E41').Value = '  -1.13%  '
$ws.Range('E42').Value = '  +19.09%  '
$ws.Range('D43').Value = "'8.98"
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('E45').Value = '  +6.19%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = "'2.54"
$ws.Range('E46').Value = '  +13.67%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').Value = "'1.27"
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = "'102.27"
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = "'1.19"
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('D50').Value = "'4.50"
$ws.Range('E50').Value = '  +0.43%  '
$ws.Range('D51').Value = "'55.54"
$ws.Range('E51').Value = '  +9.08%  '
